# Insert a new data row at row 253 (pushes existing rows 253-347 down to 254-348)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row above the current row 253; this shifts rows 253:347
# down to 254:348 (carrying their values/styles with them) exactly like a
# native Excel "Insert Row" operation.
$ws.Rows.Item(253).Insert()

# Populate the newly inserted row 253 with the new record's data.
$ws.Range("A253").Value = 10
$ws.Range("B253").Value = "Vega Modelo de Temuco"
$ws.Range("C253").Value = "La Araucanía"
$ws.Range("D253").Value = 44809
$ws.Range("D253").NumberFormat = $ws.Range("D254").NumberFormat
$ws.Range("E253").Value = 9
$ws.Range("F253").Value = 100114013
$ws.Range("G253").Value = "Zanahoria"
$ws.Range("H253").Value = "Sin especificar"
$ws.Range("I253").Value = "Primera"
$ws.Range("J253").Value = 100
$ws.Range("K253").Value = 12000
$ws.Range("L253").Value = 12000
$ws.Range("M253").Value = 12000
$ws.Range("N253").Value = "$/saco 25 kilos"
$ws.Range("O253").Value = "Región de La Araucanía"
$ws.Range("P253").Value = 480
$ws.Range("Q253").Value = 25
$ws.Range("R253").Value = "Hortaliza"
